# Generate Report for handback
# Adds a new handed-back file (3656238c-f3d4-4537-84a6-7194e5e39019) as row 4
# on the Overview, zh-cn and de-de sheets, mirroring the existing
# "b764274f-2cd2-4025-82df-3bb822e4cb38" ("in sync") entry.

$wb = $excel.ActiveWorkbook

$uuidMd      = "3656238c-f3d4-4537-84a6-7194e5e39019.md"
$zhXlf       = "3656238c-f3d4-4537-84a6-7194e5e39019.34f0377ea24abce1bbaeb5dd16f98281fb1a24c1.zh-cn.xlf"
$deXlf       = "3656238c-f3d4-4537-84a6-7194e5e39019.34f0377ea24abce1bbaeb5dd16f98281fb1a24c1.de-de.xlf"
$status      = "Handed back: in sync with en-US"
$reason      = "Include"

$zhHandoffDt  = "2016-01-19 07:39:07"
$zhHandbackDt = "2016-01-19 07:39:50"
$deHandoffDt  = "2016-01-19 07:39:18"
$deHandbackDt = "2016-01-19 07:40:08"

$overviewMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0/e2e/3656238c-f3d4-4537-84a6-7194e5e39019.md"

$zhHandoffXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1e2d3c4b5a6978869504132a1b2c3d4e5f6a7b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/3656238c-f3d4-4537-84a6-7194e5e39019.34f0377ea24abce1bbaeb5dd16f98281fb1a24c1.zh-cn.xlf"
$zhMdUrl         = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c9d8e7f6a5b4c3d2e1f0a9b8c7d6e5f4a3b2c1d0/e2e/3656238c-f3d4-4537-84a6-7194e5e39019.md"
$zhHandbackXlfUrl= "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9f8e7d6c5b4a39281706f5e4d3c2b1a09f8e7d6c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/3656238c-f3d4-4537-84a6-7194e5e39019.34f0377ea24abce1bbaeb5dd16f98281fb1a24c1.zh-cn.xlf"

$deHandoffXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a2b3c4d5e6f708192a3b4c5d6e7f8091a2b3c4d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/3656238c-f3d4-4537-84a6-7194e5e39019.34f0377ea24abce1bbaeb5dd16f98281fb1a24c1.de-de.xlf"
$deMdUrl         = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5d6e7f8091a2b3c4d5e6f708192a3b4c5d6e7f80/e2e/3656238c-f3d4-4537-84a6-7194e5e39019.md"
$deHandbackXlfUrl= "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8091a2b3c4d5e6f708192a3b4c5d6e7f8091a2b3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/3656238c-f3d4-4537-84a6-7194e5e39019.34f0377ea24abce1bbaeb5dd16f98281fb1a24c1.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet: row 4 = File Name | zh-cn status | de-de status
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4, 1).Value = $uuidMd
$wsOverview.Cells.Item(4, 2).Value = $status
$wsOverview.Cells.Item(4, 3).Value = $status

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(4, 1), $overviewMdUrl, "", "", $uuidMd)

# ---------------------------------------------------------------------------
# zh-cn sheet: row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(4, 1).Value = $uuidMd
$wsZh.Cells.Item(4, 2).Value = $status
$wsZh.Cells.Item(4, 3).Value = $zhXlf
$wsZh.Cells.Item(4, 4).Value = $zhHandoffDt
$wsZh.Cells.Item(4, 5).Value = $uuidMd
$wsZh.Cells.Item(4, 6).Value = $zhXlf
$wsZh.Cells.Item(4, 7).Value = $zhHandbackDt
$wsZh.Cells.Item(4, 8).Value = $reason

$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 1), $zhMdUrl, "", "", $uuidMd)
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 3), $zhHandoffXlfUrl, "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 5), $zhMdUrl, "", "", $uuidMd)
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 6), $zhHandbackXlfUrl, "", "", $zhXlf)

# ---------------------------------------------------------------------------
# de-de sheet: row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(4, 1).Value = $uuidMd
$wsDe.Cells.Item(4, 2).Value = $status
$wsDe.Cells.Item(4, 3).Value = $deXlf
$wsDe.Cells.Item(4, 4).Value = $deHandoffDt
$wsDe.Cells.Item(4, 5).Value = $uuidMd
$wsDe.Cells.Item(4, 6).Value = $deXlf
$wsDe.Cells.Item(4, 7).Value = $deHandbackDt
$wsDe.Cells.Item(4, 8).Value = $reason

$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 1), $deMdUrl, "", "", $uuidMd)
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 3), $deHandoffXlfUrl, "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 5), $deMdUrl, "", "", $uuidMd)
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 6), $deHandbackXlfUrl, "", "", $deXlf)
